# Scheduled market-data refresh: update price/profit columns across sheets
# per Bahamut_Profits data pull. Values sourced from the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 226913.2
$ws.Range("J17").Value = 226913.2
$ws.Range("L17").Value = 680739.6000000001
$ws.Range("N17").Value = -681075.6000000001
$ws.Range("H19").Value = 1803.4615
$ws.Range("I19").Value = 1407
$ws.Range("J19").Value = 2199.923
$ws.Range("K19").Value = 1407
$ws.Range("L19").Value = 2199.923
$ws.Range("M19").Value = -1232
$ws.Range("N19").Value = -2549.923
$ws.Range("H53").Value = 319.2857
$ws.Range("I53").Value = 95.42856999999999
$ws.Range("J53").Value = 543.1429000000001
$ws.Range("K53").Value = 95.42856999999999
$ws.Range("L53").Value = 543.1429000000001
$ws.Range("M53").Value = 541.57143
$ws.Range("N53").Value = -1817.1429
$ws.Range("H123").Value = 33363.637
$ws.Range("J123").Value = 33363.637
$ws.Range("L123").Value = 33363.637
$ws.Range("N123").Value = -43163.637
$ws.Range("H132").Value = 2569.4482
$ws.Range("I132").Value = 2493.2363
$ws.Range("J132").Value = 3966.6667
$ws.Range("K132").Value = 7479.7089
$ws.Range("L132").Value = 11900.0001
$ws.Range("M132").Value = -4949.7089
$ws.Range("N132").Value = -16960.0001
$ws.Range("H137").Value = 803.9524
$ws.Range("I137").Value = 803.9524
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2411.8572
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 1314.6061
$ws.Range("I138").Value = 695.0714
$ws.Range("J138").Value = 2121.442
$ws.Range("K138").Value = 2085.2142
$ws.Range("L138").Value = 6364.326
$ws.Range("M138").Value = 3054.7858
$ws.Range("N138").Value = -16644.326

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4409631
$ws.Range("I32").Value = 5673618
$ws.Range("J32").Value = 22853.646
$ws.Range("K32").Value = 5673618
$ws.Range("L32").Value = 22853.646
$ws.Range("M32").Value = -5673331
$ws.Range("N32").Value = -23427.646
$ws.Range("H37").Value = 7833.1665
$ws.Range("I37").Value = 7000
$ws.Range("J37").Value = 7999.8
$ws.Range("K37").Value = 7000
$ws.Range("L37").Value = 7999.8
$ws.Range("M37").Value = -6727
$ws.Range("N37").Value = -8545.799999999999
$ws.Range("H44").Value = 19762.25
$ws.Range("J44").Value = 19762.25
$ws.Range("L44").Value = 19762.25
$ws.Range("N44").Value = -20738.25
$ws.Range("H45").Value = 1797.1428
$ws.Range("I45").Value = 1263.3334
$ws.Range("K45").Value = 1263.3334
$ws.Range("M45").Value = -886.3334
$ws.Range("H55").Value = 149551.25
$ws.Range("J55").Value = 149551.25
$ws.Range("L55").Value = 149551.25
$ws.Range("N55").Value = -150181.25
$ws.Range("H110").Value = 533.36
$ws.Range("I110").Value = 465
$ws.Range("K110").Value = 465
$ws.Range("M110").Value = 1580
$ws.Range("H122").Value = 867.0714
$ws.Range("I122").Value = 879.0833
$ws.Range("J122").Value = 795
$ws.Range("K122").Value = 2637.2499
$ws.Range("L122").Value = 2385
$ws.Range("M122").Value = -187.2498999999998
$ws.Range("N122").Value = -7285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H99").Value = 66668730
$ws.Range("I99").Value = 90911210
$ws.Range("J99").Value = 1905.5
$ws.Range("K99").Value = 90911210
$ws.Range("L99").Value = 1905.5
$ws.Range("M99").Value = -90909712
$ws.Range("N99").Value = -4901.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2696.0303
$ws.Range("I31").Value = 2740.9355
$ws.Range("K31").Value = 2740.9355
$ws.Range("M31").Value = -2445.9355
$ws.Range("H34").Value = 2696.0303
$ws.Range("I34").Value = 2740.9355
$ws.Range("K34").Value = 2740.9355
$ws.Range("M34").Value = -2538.9355
$ws.Range("H62").Value = 6251.4287
$ws.Range("I62").Value = 5960
$ws.Range("K62").Value = 5960
$ws.Range("M62").Value = -5336
$ws.Range("H65").Value = 6251.4287
$ws.Range("I65").Value = 5960
$ws.Range("K65").Value = 29800
$ws.Range("M65").Value = -26680
$ws.Range("H94").Value = 2525.3447
$ws.Range("I94").Value = 4253
$ws.Range("J94").Value = 2248.92
$ws.Range("K94").Value = 4253
$ws.Range("L94").Value = 2248.92
$ws.Range("M94").Value = -3802
$ws.Range("N94").Value = -3150.92
$ws.Range("H99").Value = 2446.7576
$ws.Range("I99").Value = 2202.7778
$ws.Range("J99").Value = 2739.5334
$ws.Range("K99").Value = 2202.7778
$ws.Range("L99").Value = 2739.5334
$ws.Range("M99").Value = -704.7777999999998
$ws.Range("N99").Value = -5735.5334
$ws.Range("H126").Value = 2446.7576
$ws.Range("I126").Value = 2202.7778
$ws.Range("J126").Value = 2739.5334
$ws.Range("K126").Value = 6608.3334
$ws.Range("L126").Value = 8218.600199999999
$ws.Range("M126").Value = -4138.3334
$ws.Range("N126").Value = -13158.6002
$ws.Range("H134").Value = 1182.5714
$ws.Range("I134").Value = 1031.4166
$ws.Range("J134").Value = 1666.2667
$ws.Range("K134").Value = 3094.2498
$ws.Range("L134").Value = 4998.800099999999
$ws.Range("M134").Value = -559.2498000000001
$ws.Range("N134").Value = -10068.8001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2108.56
$ws.Range("J129").Value = 2623.6
$ws.Range("L129").Value = 7870.799999999999
$ws.Range("N129").Value = -17870.8
$ws.Range("H131").Value = 794.84
$ws.Range("J131").Value = 834.5165
$ws.Range("L131").Value = 2503.5495
$ws.Range("N131").Value = -12583.5495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4225.55
$ws.Range("I70").Value = 3942.6667
$ws.Range("K70").Value = 3942.6667
$ws.Range("M70").Value = -3672.6667
$ws.Range("H73").Value = 4225.55
$ws.Range("I73").Value = 3942.6667
$ws.Range("K73").Value = 3942.6667
$ws.Range("M73").Value = -3006.6667
$ws.Range("H80").Value = 2934.3
$ws.Range("I80").Value = 3305
$ws.Range("J80").Value = 1451.5
$ws.Range("K80").Value = 3305
$ws.Range("L80").Value = 1451.5
$ws.Range("M80").Value = -2307
$ws.Range("N80").Value = -3447.5
$ws.Range("H83").Value = 2934.3
$ws.Range("I83").Value = 3305
$ws.Range("J83").Value = 1451.5
$ws.Range("K83").Value = 16525
$ws.Range("L83").Value = 7257.5
$ws.Range("M83").Value = -11533
$ws.Range("N83").Value = -17241.5
$ws.Range("H126").Value = 2728.9333
$ws.Range("I126").Value = 2940.375
$ws.Range("J126").Value = 2487.2856
$ws.Range("K126").Value = 8821.125
$ws.Range("L126").Value = 7461.8568
$ws.Range("M126").Value = -6351.125
$ws.Range("N126").Value = -12401.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3459.889
$ws.Range("I96").Value = 2509
$ws.Range("J96").Value = 4954.143
$ws.Range("K96").Value = 2509
$ws.Range("L96").Value = 4954.143
$ws.Range("M96").Value = -1136
$ws.Range("N96").Value = -7700.143
$ws.Range("H122").Value = 1208.5385
$ws.Range("I122").Value = 1000.8571
$ws.Range("J122").Value = 1450.8334
$ws.Range("K122").Value = 3002.5713
$ws.Range("L122").Value = 4352.5002
$ws.Range("M122").Value = -552.5712999999996
$ws.Range("N122").Value = -9252.5002
$ws.Range("H132").Value = 1859.3334
$ws.Range("I132").Value = 1052.8
$ws.Range("J132").Value = 3100.1538
$ws.Range("K132").Value = 3158.4
$ws.Range("L132").Value = 9300.4614
$ws.Range("M132").Value = -628.3999999999996
$ws.Range("N132").Value = -14360.4614
